# Fruta / hortaliza, semanal
# Insert a new weekly record at row 595 (shifting existing rows 595-643 down
# to 596-644) and populate the new row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 595; everything below shifts
# down by one (old 595 -> 596, ..., old 643 -> 644).
$ws.Rows.Item(595).Insert()

$row = 595

$ws.Cells.Item($row, 1).Value2 = 9
$ws.Cells.Item($row, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value2 = "Metropolitana"
$ws.Cells.Item($row, 4).Value2 = 45223
$ws.Cells.Item($row, 5).Value2 = 13
$ws.Cells.Item($row, 6).Value2 = 100112039
$ws.Cells.Item($row, 7).Value2 = "Ciboulette"
$ws.Cells.Item($row, 8).Value2 = "Sin especificar"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 340
$ws.Cells.Item($row, 11).Value2 = 1200
$ws.Cells.Item($row, 12).Value2 = 1400
$ws.Cells.Item($row, 13).Value2 = 1300
$ws.Cells.Item($row, 14).Value2 = "$/docena de atados"
$ws.Cells.Item($row, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value2 = 433
$ws.Cells.Item($row, 17).Value2 = 3
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
